$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (17) of data to the rental contracts table, mirroring the
# existing rows: ID, address, area, purchase price, property, deposit,
# tenant, expiry date.
$ws.Range("A17").Value = "A04"
$ws.Range("B17").Value = "테스트"
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = "소망"
$ws.Range("F17").Value = 580000000
$ws.Range("G17").Value = "명민재"
$ws.Range("H17").Value = (Get-Date -Year 2025 -Month 8 -Day 29 -Hour 0 -Minute 0 -Second 0)

# Copy formatting from the row above (row 16) so the new row matches the
# existing table styling.
$ws.Range("A16:H16").Copy() | Out-Null
$ws.Range("A17:H17").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# Update the selected cell to match the post-edit selection.
$ws.Range("H18").Select() | Out-Null
